$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (rows 289, 322, 327, 328)
$ws.Range("B289").Value = 32957600000
$ws.Range("D289").Value = 46550282485.87571

$ws.Range("B322").Value = 36383400000
$ws.Range("D322").Value = 51403503814.6369

$ws.Range("B327").Value = 37472200000
$ws.Range("D327").Value = 52941791466.51597

$ws.Range("B328").Value = 37492700000
$ws.Range("D328").Value = 52970754450.40972

# Append new rows 352-353, matching the style of column A (date style "2")
$ws.Range("A352").Value = 44986
$ws.Range("B352").Value = 41710700000
$ws.Range("C352").Value = 1.409443269908386
$ws.Range("D352").Value = 58788865398.16772

$ws.Range("A353").Value = 45017
$ws.Range("B353").Value = 41839100000
$ws.Range("C353").Value = 1.410835214446953
$ws.Range("D353").Value = 59028075620.7675

# Copy style from A351 to A352:A353 so the date format / alignment match
$ws.Range("A351").Copy()
$ws.Range("A352:A353").PasteSpecial(-4122)  # xlPasteFormats
